$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing columns (B:H) for rows 2 and 3 with refreshed data ---
$ws.Range("B2").Value = 206
$ws.Range("C2").Value = 5354513
$ws.Range("D2").Value = 484651.17839999998
$ws.Range("E2").Value = 51257
$ws.Range("F2").Value = 21946.979000000003
$ws.Range("G2").Value = 68645276
$ws.Range("H2").Value = 142.53555555555556

$ws.Range("B3").Value = 163
$ws.Range("C3").Value = 4739700
$ws.Range("D3").Value = 241356.24260000003
$ws.Range("E3").Value = 11750
$ws.Range("F3").Value = 35644.804000000004
$ws.Range("G3").Value = 47481082
$ws.Range("H3").Value = 37.783611111111114

# --- New percentage-breakdown columns (I:O) ---
$ws.Range("I1").Value = "ClassNumberVideos_percentage"
$ws.Range("J1").Value = "ClassViews_percentage"
$ws.Range("K1").Value = "ClassWatchTime_hr_percentage"
$ws.Range("L1").Value = "ClassSubscribers_percentage"
$ws.Range("M1").Value = "ClassRevenue_USD_percentage"
$ws.Range("N1").Value = "ClassImpressions_percentage"
$ws.Range("O1").Value = "ClassVideoDuration_hr_percentage"

$ws.Range("I2").Value = 55.826558265582662
$ws.Range("J2").Value = 53.04537362149977
$ws.Range("K2").Value = 66.755678300428826
$ws.Range("L2").Value = 81.351278429380855
$ws.Range("M2").Value = 38.107830417405204
$ws.Range("N2").Value = 59.112571152881586
$ws.Range("O2").Value = 79.04625902527772

$ws.Range("I3").Value = 44.173441734417345
$ws.Range("J3").Value = 46.954626378500237
$ws.Range("K3").Value = 33.244321699571174
$ws.Range("L3").Value = 18.648721570619138
$ws.Range("M3").Value = 61.892169582594789
$ws.Range("N3").Value = 40.887428847118414
$ws.Range("O3").Value = 20.953740974722294

# --- Column widths for the new columns (values chosen so the engine's
#     ColumnWidth -> stored-width rounding lands on the target widths) ---
$ws.Columns.Item(9).ColumnWidth = 29.166666666666668
$ws.Columns.Item(10).ColumnWidth = 21.166666666666668
$ws.Columns.Item(11).ColumnWidth = 28.666666666666668
$ws.Columns.Item(12).ColumnWidth = 25.833333333333332
$ws.Columns.Item(13).ColumnWidth = 28.166666666666668
$ws.Columns.Item(14).ColumnWidth = 26.5
$ws.Columns.Item(15).ColumnWidth = 31.666666666666668
